$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Equip-001
$ws.Range("B2").Value = "DESC-8"
$ws.Range("C2").Value = "PT-4"
$ws.Range("D2").Value = "PDT-5"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 8

# Row 3 - Equip-002
$ws.Range("B3").Value = "DESC-5"
$ws.Range("C3").Value = "PT-2"
$ws.Range("D3").Value = "PDT-4"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 7

# Row 4 - Equip-003
$ws.Range("B4").Value = "DESC-6"
$ws.Range("C4").Value = "PT-6"
$ws.Range("D4").Value = "PDT-5"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 3

# Row 5 - Equip-004
$ws.Range("B5").Value = "DESC-7"
$ws.Range("C5").Value = "PT-5"
$ws.Range("D5").Value = "PDT-4"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = 9
